$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "21.026.55"
$ws.Range("E2").Value = "  -4.71%  "

$ws.Range("D3").Value = "1.501.75"
$ws.Range("E3").Value = "  -3.48%  "

$ws.Range("E4").Value = "  +0.70%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.006"
$ws.Range("E5").Value = "  +0.54%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "283.47"
$ws.Range("E6").Value = "  -2.46%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3815"
$ws.Range("E7").Value = "  -3.55%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3116"
$ws.Range("E8").Value = "  -3.39%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "42.71"
$ws.Range("E9").Value = "  -2.54%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06909"
$ws.Range("E10").Value = "  -4.89%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.032"
$ws.Range("E11").Value = "  -4.38%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.007"
$ws.Range("E12").Value = "  +0.76%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.570"
$ws.Range("E13").Value = "  -2.28%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.71"
$ws.Range("E14").Value = "  -5.94%  "

$ws.Range("D15").Value = "1.508.99"
$ws.Range("E15").Value = "  -3.00%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.355"
$ws.Range("E16").Value = "  -4.27%  "

$ws.Range("E17").Value = "  -6.06%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06553"
$ws.Range("E18").Value = "  -0.62%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "81.79"
$ws.Range("E19").Value = "  -2.18%  "

$ws.Range("E20").Value = "  +0.57%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.955"
$ws.Range("E21").Value = "  -5.29%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.03"
$ws.Range("E22").Value = "  -3.34%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.84"
$ws.Range("E23").Value = "  -4.34%  "

$ws.Range("E24").Value = "  -1.30%  "

$ws.Range("D25").Value = "21.027.05"
$ws.Range("E25").Value = "  -4.74%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.317"
$ws.Range("E26").Value = "  -4.59%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "146.94"
$ws.Range("E27").Value = "  -1.30%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.93"
$ws.Range("E28").Value = "  -4.06%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.800"
$ws.Range("E29").Value = "  -1.74%  "

$ws.Range("D30").Value = "1.676.61"
$ws.Range("E30").Value = "  -3.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "114.14"
$ws.Range("E31").Value = "  -4.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.828"
$ws.Range("E32").Value = "  -0.35%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9479"
$ws.Range("E33").Value = "  -3.02%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07899"
$ws.Range("E34").Value = "  -5.17%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "8.403"
$ws.Range("E35").Value = "  -7.92%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").Value = "5.047"
$ws.Range("E36").Value = "  -1.55%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("B37").Value = "WEMIXTOKEN"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").Value = "1.473"
$ws.Range("E37").Value = "  -8.16%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "11.12"
$ws.Range("E38").Value = "  +3.74%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05760"
$ws.Range("E39").Value = "  -4.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.006"
$ws.Range("E40").Value = "  +0.29%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02116"
$ws.Range("E41").Value = "  -6.76%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.157"
$ws.Range("E42").Value = "  -4.32%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1964"
$ws.Range("E43").Value = "  -3.62%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5589"
$ws.Range("E44").Value = "  -4.12%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.81"
$ws.Range("E45").Value = "  -1.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.658"
$ws.Range("E46").Value = "  -2.38%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5406"
$ws.Range("E47").Value = "  -3.29%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "1.837"
$ws.Range("E48").Value = "  -3.51%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("B49").Value = "EOS"
$ws.Range("C49").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D49").Value = "1.127"
$ws.Range("E49").Value = "  -0.67%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "112.95"
$ws.Range("E50").Value = "  -4.48%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06526"
$ws.Range("E51").Value = "  -4.30%  "

